$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header: insert a new "Trening" header in column B, and move the
# "Acceleration_SMA" header to column C (reusing the same bold/bordered
# header style that was already on A1/B1).
$ws.Range("C1").Value = "Acceleration_SMA"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Trening"

# Data rows: Velocity_Bin (A), Trening (B), Acceleration_SMA (C)
$data = @(
    @("10-15", "Gra",      3.714678005377452),
    @("5-10",  "Gra",      3.347213657555126),
    @("10-15", "Duża Gra", 1.934777881311519),
    @("10-15", "Mała Gra", 3.27617809676885),
    @("5-10",  "Duża Gra", 1.862376963063365),
    @("5-10",  "Mała Gra", 2.847562606920276),
    @("10-15", "Duża Gra", 3.841688387559884),
    @("10-15", "Mała Gra", 3.389416957949544),
    @("5-10",  "Duża Gra", 3.266047984689147),
    @("5-10",  "Mała Gra", 2.928736258011598),
    @("10-15", "Duża Gra", 1.131088749854399),
    @("10-15", "Mała Gra", 3.073714245217187),
    @("5-10",  "Duża Gra", 1.606974919637044),
    @("5-10",  "Mała Gra", 2.781481812965303)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
